$d = $word.ActiveDocument

# 1) Update the NOMENCLATURE label for the "status" field
$d.Content.Find.Execute(
    "(NOMENCLATURE: ENUM-STATUS_DR)", $false, $false, $false, $false, $false,
    $true, 1, $false, "(NOMENCLATURE: HubSante.etatDemande)", 2)

# 2) Replace the description text for the "status" field.
#    Done via Range.Text (rather than Find/Replace) so that Word's
#    AutoCorrect "smart quotes" feature does not mangle the apostrophe.
$rng = $d.Content
$rng.Find.Execute("A valoriser avec la valeur ANNULEE uniquement pour signifier l'annulation d'une demande de ressources. Les autres champs de la demande sont remplis à l'identique de la demande initiale envoyée.")
$rng.Text = "A quoi ça sert d'avoir un objet demande "

# 3) Update the NOMENCLATURE label for the "convention" field
$d.Content.Find.Execute(
    "(NOMENCLATURE: CISU-CADRE_CONV)", $false, $false, $false, $false, $false,
    $true, 1, $false, "(NOMENCLATURE: HubSante.cadre)", 2)

# 4) Update the NOMENCLATURE label for the "purpose" field
$d.Content.Find.Execute(
    "(NOMENCLATURE: CISU-Code_Effet_a_obtenir)", $false, $false, $false, $false, $false,
    $true, 1, $false, "(NOMENCLATURE: HubSante.effet)", 2)

# 5) Update the NOMENCLATURE label for the "deadline" field
$d.Content.Find.Execute(
    "(NOMENCLATURE: SI-SAMU-DELAI)", $false, $false, $false, $false, $false,
    $true, 1, $false, "(NOMENCLATURE: HubSante.delai)", 2)
